# Updates the cryptos list (price / volume-1h columns) per the latest
# scrape. A handful of D-column price strings look numeric (e.g. "226.23",
# "0.0500") so we force those cells to Text format first -- otherwise Excel
# auto-coerces them to numbers and silently drops significant trailing
# zeros (e.g. "0.0500" -> 0.05), which would not match the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.614.67'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.812.60'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.23'
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '37.96'
$ws.Range("E8").Value = '  +8.56%  '
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0682'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0970'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").Value = '2.074.02'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").Value = '1.815.22'
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").Value = '34.570.20'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.71'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.76'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '0.0₃0778'
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("E21").Value = '  -1.75%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  +4.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.88'
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.90'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("E27").Value = '  +3.52%  '
$ws.Range("E28").Value = '  +1.59%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.24'
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("D35").Value = '1.367.12'
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.656'
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -3.92%  '
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("E40").Value = '  +8.79%  '
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.21'
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.941'
$ws.Range("E44").Value = '  -2.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.97'
$ws.Range("E45").Value = '  +4.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0500'
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("D47").Value = '1.974.88'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.11'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("E51").Value = '  -7.53%  '
